$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1551.8235
$ws.Cells.Item(19, 9).Value = 1554.091
$ws.Cells.Item(19, 10).Value = 1547.6666
$ws.Cells.Item(19, 11).Value = 1554.091
$ws.Cells.Item(19, 12).Value = 1547.6666
$ws.Cells.Item(19, 13).Value = -1379.091
$ws.Cells.Item(19, 14).Value = -1897.6666
$ws.Cells.Item(74, 8).Value = 5616.892
$ws.Cells.Item(74, 9).Value = 4447.7646
$ws.Cells.Item(74, 10).Value = 6610.65
$ws.Cells.Item(74, 11).Value = 4447.7646
$ws.Cells.Item(74, 12).Value = 6610.65
$ws.Cells.Item(74, 13).Value = -3511.7646
$ws.Cells.Item(74, 14).Value = -8482.65
$ws.Cells.Item(77, 8).Value = 5616.892
$ws.Cells.Item(77, 9).Value = 4447.7646
$ws.Cells.Item(77, 10).Value = 6610.65
$ws.Cells.Item(77, 11).Value = 22238.823
$ws.Cells.Item(77, 12).Value = 33053.25
$ws.Cells.Item(77, 13).Value = -17558.823
$ws.Cells.Item(77, 14).Value = -42413.25
$ws.Cells.Item(98, 8).Value = 2520.279
$ws.Cells.Item(98, 9).Value = 2071.3235
$ws.Cells.Item(98, 11).Value = 2071.3235
$ws.Cells.Item(98, 13).Value = -573.3235
$ws.Cells.Item(122, 8).Value = 2520.279
$ws.Cells.Item(122, 9).Value = 2071.3235
$ws.Cells.Item(122, 11).Value = 6213.970499999999
$ws.Cells.Item(122, 13).Value = -3763.970499999999
$ws.Cells.Item(138, 8).Value = 9858.111000000001
$ws.Cells.Item(138, 9).Value = 8149.5
$ws.Cells.Item(138, 10).Value = 10227.541
$ws.Cells.Item(138, 11).Value = 24448.5
$ws.Cells.Item(138, 12).Value = 30682.623
$ws.Cells.Item(138, 13).Value = -19308.5
$ws.Cells.Item(138, 14).Value = -40962.623
$ws.Cells.Item(141, 8).Value = 7538.2
$ws.Cells.Item(141, 9).Value = 4687.4
$ws.Cells.Item(141, 10).Value = 13239.8
$ws.Cells.Item(141, 11).Value = 14062.2
$ws.Cells.Item(141, 12).Value = 39719.39999999999
$ws.Cells.Item(141, 13).Value = -8882.199999999999
$ws.Cells.Item(141, 14).Value = -50079.39999999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1653.0264
$ws.Cells.Item(32, 9).Value = 1032.2727
$ws.Cells.Item(32, 11).Value = 1032.2727
$ws.Cells.Item(32, 13).Value = -745.2727
$ws.Cells.Item(61, 8).Value = 5959.161
$ws.Cells.Item(61, 10).Value = 9503.223
$ws.Cells.Item(61, 12).Value = 9503.223
$ws.Cells.Item(61, 14).Value = -9927.223
$ws.Cells.Item(74, 8).Value = 2894
$ws.Cells.Item(74, 9).Value = 1778.5
$ws.Cells.Item(74, 10).Value = 4381.3335
$ws.Cells.Item(74, 11).Value = 1778.5
$ws.Cells.Item(74, 12).Value = 4381.3335
$ws.Cells.Item(74, 13).Value = -904.5
$ws.Cells.Item(74, 14).Value = -6129.3335
$ws.Cells.Item(77, 8).Value = 2894
$ws.Cells.Item(77, 9).Value = 1778.5
$ws.Cells.Item(77, 10).Value = 4381.3335
$ws.Cells.Item(77, 11).Value = 8892.5
$ws.Cells.Item(77, 12).Value = 21906.6675
$ws.Cells.Item(77, 13).Value = -4524.5
$ws.Cells.Item(77, 14).Value = -30642.6675
$ws.Cells.Item(110, 8).Value = 59300.844
$ws.Cells.Item(110, 9).Value = 65992.766
$ws.Cells.Item(110, 10).Value = 2419.5
$ws.Cells.Item(110, 11).Value = 65992.766
$ws.Cells.Item(110, 12).Value = 2419.5
$ws.Cells.Item(110, 13).Value = -63947.766
$ws.Cells.Item(110, 14).Value = -6509.5
$ws.Cells.Item(136, 8).Value = 5959.161
$ws.Cells.Item(136, 10).Value = 9503.223
$ws.Cells.Item(136, 12).Value = 28509.669
$ws.Cells.Item(136, 14).Value = -33609.669

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3580.3333
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 3580.3333
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 3580.3333
$ws.Cells.Item(16, 13).ClearContents()
$ws.Cells.Item(16, 14).Value = -4154.3333
$ws.Cells.Item(31, 8).Value = 3669.8774
$ws.Cells.Item(31, 9).Value = 3389.8
$ws.Cells.Item(31, 10).Value = 3741.6924
$ws.Cells.Item(31, 11).Value = 3389.8
$ws.Cells.Item(31, 12).Value = 3741.6924
$ws.Cells.Item(31, 13).Value = -3094.8
$ws.Cells.Item(31, 14).Value = -4331.6924
$ws.Cells.Item(34, 8).Value = 3669.8774
$ws.Cells.Item(34, 9).Value = 3389.8
$ws.Cells.Item(34, 10).Value = 3741.6924
$ws.Cells.Item(34, 11).Value = 3389.8
$ws.Cells.Item(34, 12).Value = 3741.6924
$ws.Cells.Item(34, 13).Value = -3187.8
$ws.Cells.Item(34, 14).Value = -4145.6924
$ws.Cells.Item(109, 8).Value = 62142.5
$ws.Cells.Item(109, 10).Value = 84285
$ws.Cells.Item(109, 12).Value = 84285
$ws.Cells.Item(109, 14).Value = -86365
$ws.Cells.Item(113, 8).Value = 3580.3333
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 3580.3333
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 3580.3333
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).Value = -7920.3333
$ws.Cells.Item(134, 8).Value = 4021.8857
$ws.Cells.Item(134, 9).Value = 3568.6667
$ws.Cells.Item(134, 10).Value = 11500
$ws.Cells.Item(134, 11).Value = 10706.0001
$ws.Cells.Item(134, 12).Value = 34500
$ws.Cells.Item(134, 13).Value = -8171.000100000001
$ws.Cells.Item(134, 14).Value = -39570

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(109, 8).Value = 3583.7778
$ws.Cells.Item(109, 9).Value = 2188.5
$ws.Cells.Item(109, 10).Value = 6374.3335
$ws.Cells.Item(109, 11).Value = 6565.5
$ws.Cells.Item(109, 12).Value = 19123.0005
$ws.Cells.Item(109, 13).Value = -5525.5
$ws.Cells.Item(109, 14).Value = -21203.0005

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 4026.6
$ws.Cells.Item(126, 9).Value = 2295.125
$ws.Cells.Item(126, 10).Value = 6005.4287
$ws.Cells.Item(126, 11).Value = 6885.375
$ws.Cells.Item(126, 12).Value = 18016.2861
$ws.Cells.Item(126, 13).Value = -4415.375
$ws.Cells.Item(126, 14).Value = -22956.2861

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1000
$ws.Cells.Item(61, 9).Value = 1000
$ws.Cells.Item(61, 11).Value = 1000
$ws.Cells.Item(61, 13).Value = -798
$ws.Cells.Item(96, 8).Value = 83598
$ws.Cells.Item(96, 10).Value = 83598
$ws.Cells.Item(96, 12).Value = 83598
$ws.Cells.Item(96, 14).Value = -89090
$ws.Cells.Item(109, 8).Value = 69235.664
$ws.Cells.Item(109, 10).Value = 69235.664
$ws.Cells.Item(109, 12).Value = 69235.664
$ws.Cells.Item(109, 14).Value = -72009.664
$ws.Cells.Item(113, 8).Value = 1000
$ws.Cells.Item(113, 9).Value = 1000
$ws.Cells.Item(113, 11).Value = 1000
$ws.Cells.Item(113, 13).Value = 1170

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 58249.75
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 13).ClearContents()
$ws.Cells.Item(73, 8).Value = 58249.75
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 17270
$ws.Cells.Item(132, 9).Value = 50000
$ws.Cells.Item(132, 10).Value = 6360
$ws.Cells.Item(132, 11).Value = 150000
$ws.Cells.Item(132, 12).Value = 19080
$ws.Cells.Item(132, 13).Value = -147470
$ws.Cells.Item(132, 14).Value = -24140
$ws.Cells.Item(136, 8).Value = 7038.8125
$ws.Cells.Item(136, 9).Value = 6870.9023
$ws.Cells.Item(136, 10).Value = 8022.2856
$ws.Cells.Item(136, 11).Value = 20612.7069
$ws.Cells.Item(136, 12).Value = 24066.8568
$ws.Cells.Item(136, 13).Value = -18062.7069
$ws.Cells.Item(136, 14).Value = -29166.8568

Write-Output "Applied 163 sets and 4 clears"
